# ccm_RLC.xlsx — "Add files via upload" edit
#
# 1) L_map: bump group-index column B for rows 2-10 (0/1 -> 1/2)
# 2) Sheet re-numbering: SS0 -> SS1, (old) SS1 -> SS2, plus a brand new SS3
#    tab that is a copy of SS2's content.
# 3) Minor view/selection bookkeeping to mirror the saved workbook state.

$wb = $excel.ActiveWorkbook

# --- 1) L_map data edits -------------------------------------------------
$ws = $wb.Worksheets.Item("L_map")
$ws.Range("B2:B5").Value = 1
$ws.Range("B6:B10").Value = 2
$ws.Range("B11").Select()

# --- 2) Rename sheets (do the last one first so the names never collide) -
$ssOld1 = $wb.Worksheets.Item("SS1")     # will become SS2
$ssOld1.Name = "SS2"

$ssOld0 = $wb.Worksheets.Item("SS0")     # will become SS1
$ssOld0.Name = "SS1"

# Re-normalise A5's redundant "border + fill" format down to the plain
# bordered style already used elsewhere on the sheet (B5/etc use it) -
# matches the de-duped style table Excel wrote back on save.
$ss2 = $wb.Worksheets.Item("SS2")
$ss2.Range("A5").Borders.Item(7).LineStyle = 1
$ss2.Range("A5").Borders.Item(9).LineStyle = 1

# --- 3) Duplicate SS2 to create the new SS3 tab ---------------------------
$ss2.Copy($null, $ss2)
$ss3 = $wb.Worksheets.Item("SS2 (2)")
$ss3.Name = "SS3"
$ss3.Range("A1:D5").Select()

# --- 4) Restore cosmetic selections on the other tabs ---------------------
$ss1 = $wb.Worksheets.Item("SS1")
$ss1.Range("E12").Select()

$ss2.Range("C3").Select()

$ss1.Activate()
$ws.Activate()
